$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1436797380447386
$ws.Range("D2").Value = 0.7403357923030853
$ws.Range("E2").Value = -1.325827866792679
$ws.Range("F2").Value = 0.0355829000473022
$ws.Range("G2").Value = 0.3640756905078888
$ws.Range("H2").Value = -0.06383541971445079
$ws.Range("C3").Value = -1.716686248779299
$ws.Range("D3").Value = 1.053612291812897
$ws.Range("E3").Value = 0.4070562124252333
$ws.Range("F3").Value = 0.5752823352813721
$ws.Range("G3").Value = 0.5053382515907288
$ws.Range("H3").Value = -0.007177666760981
$ws.Range("C4").Value = -3.969705402851109
$ws.Range("D4").Value = 0.09983259439468162
$ws.Range("E4").Value = 0.9274015724658967
$ws.Range("F4").Value = 0.5288565754890442
$ws.Range("G4").Value = 0.7619016766548157
$ws.Range("H4").Value = 0.2580905556678772
$ws.Range("C5").Value = -5.239097833633425
$ws.Range("D5").Value = -0.5959589481353762
$ws.Range("E5").Value = -0.1143757104873694
$ws.Range("F5").Value = -0.7357872128486633
$ws.Range("G5").Value = -0.0797179117798805
$ws.Range("H5").Value = 0.8136724829673767
$ws.Range("C6").Value = 2.036354780197156
$ws.Range("D6").Value = -0.2977316975593562
$ws.Range("E6").Value = -1.156379550695419
$ws.Range("F6").Value = -1.146135926246643
$ws.Range("G6").Value = -0.6188064813613892
$ws.Range("H6").Value = -0.1640173196792602
$ws.Range("C7").Value = 5.741946458816498
$ws.Range("D7").Value = -4.677600264549274
$ws.Range("E7").Value = -8.871290028095274
$ws.Range("F7").Value = 0.8213083148002625
$ws.Range("G7").Value = -3.046079635620117
$ws.Range("H7").Value = 1.09803032875061
$ws.Range("C8").Value = -18.10554087162019
$ws.Range("D8").Value = -17.49406802654266
$ws.Range("E8").Value = -29.53683829307556
$ws.Range("F8").Value = 1.518305540084839
$ws.Range("G8").Value = -0.5958990454673767
$ws.Range("H8").Value = 0.1006400510668754
$ws.Range("C9").Value = 1.688319206237781
$ws.Range("D9").Value = -2.673514366149895
$ws.Range("E9").Value = -12.34325218200682
$ws.Range("F9").Value = -0.9292787313461304
$ws.Range("G9").Value = 0.0415388382971286
$ws.Range("H9").Value = -2.335643291473389
$ws.Range("C10").Value = -3.507262408733383
$ws.Range("D10").Value = 0.2684899270534595
$ws.Range("E10").Value = -2.699394106864903
$ws.Range("F10").Value = 0.5500841736793518
$ws.Range("G10").Value = 1.703092336654663
$ws.Range("H10").Value = -0.4928155243396759
$ws.Range("C11").Value = -0.7152169942855791
$ws.Range("D11").Value = 0.4368197321891781
$ws.Range("E11").Value = 0.5279676914215112
$ws.Range("F11").Value = 0.2722931802272796
$ws.Range("G11").Value = -0.3020728528499603
$ws.Range("H11").Value = 0.2756529450416565
$ws.Range("C12").Value = 0.77803122997284
$ws.Range("D12").Value = 0.2878375947475432
$ws.Range("E12").Value = -0.8067402243614243
$ws.Range("F12").Value = -0.6982190012931824
$ws.Range("G12").Value = 0.4492913782596588
$ws.Range("H12").Value = -0.5900958180427551
$ws.Range("C13").Value = -0.02435183525085627
$ws.Range("D13").Value = -0.3931519985198989
$ws.Range("E13").Value = -2.55875074863434
$ws.Range("F13").Value = -0.4308127164840698
$ws.Range("G13").Value = 1.22447943687439
$ws.Range("H13").Value = -0.319024384021759
$ws.Range("C14").Value = -0.2728092074394176
$ws.Range("D14").Value = -0.8767854124307625
$ws.Range("E14").Value = -2.626779749989506
$ws.Range("F14").Value = -0.0684169083833694
$ws.Range("G14").Value = 0.4977024495601654
$ws.Range("H14").Value = -0.0363464802503585
$ws.Range("C15").Value = 1.066039085388185
$ws.Range("D15").Value = -0.5432969331741325
$ws.Range("E15").Value = -1.81298840045929
$ws.Range("F15").Value = -0.1298088580369949
$ws.Range("G15").Value = 0.0172569435089826
$ws.Range("H15").Value = 0.1985312104225158
$ws.Range("C16").Value = 1.15899240970611
$ws.Range("D16").Value = -0.2723855525255186
$ws.Range("E16").Value = -1.873978555202483
$ws.Range("F16").Value = 0.0598647929728031
$ws.Range("G16").Value = 0.5253441333770752
$ws.Range("H16").Value = -0.0415388382971286
$ws.Range("C17").Value = -0.3004360198974618
$ws.Range("D17").Value = 0.163273096084592
$ws.Range("E17").Value = -1.427715420722963
$ws.Range("F17").Value = -0.1009454801678657
$ws.Range("G17").Value = 0.2724458873271942
$ws.Range("H17").Value = -0.4051563739776611
$ws.Range("C18").Value = -0.5053797960281373
$ws.Range("D18").Value = -0.4524855315685276
$ws.Range("E18").Value = -1.996987149119378
$ws.Range("F18").Value = -0.0010690141934901
$ws.Range("G18").Value = -0.1505782902240753
$ws.Range("H18").Value = 0.0145080499351024
$ws.Range("C19").Value = 0.1984210014343271
$ws.Range("D19").Value = -0.3136261999607083
$ws.Range("E19").Value = -1.604727536439895
$ws.Range("F19").Value = -0.1108720451593399
$ws.Range("G19").Value = 0.0652098655700683
$ws.Range("H19").Value = -0.131183311343193
$ws.Range("C20").Value = 0.1348390579223627
$ws.Range("D20").Value = -0.3085210472345353
$ws.Range("E20").Value = -1.704802349209786
$ws.Range("F20").Value = 0.0563523173332214
$ws.Range("G20").Value = 0.0316122770309448
$ws.Range("H20").Value = 0.1798998117446899
$ws.Range("C21").Value = -0.11956262588501
$ws.Range("D21").Value = -0.3192775845527647
$ws.Range("E21").Value = -1.924065947532654
$ws.Range("F21").Value = -0.0160352122038602
$ws.Range("G21").Value = 0.0537561401724815
$ws.Range("H21").Value = -0.0355829000473022
